# Add 15 new alumni rows (rows 10-24) to Sheet1, mirroring the formatting
# of the existing data rows but with wrap-text enabled on column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Vivek Gusain",
    "Marmik Sharma",
    "Mayank Singh",
    "Mukul C. Mahadik",
    "Sumit Sharma",
    "Abhiraj Singh Rathore",
    "Amol Bobade",
    "Divyanshu Bhaik",
    "Kunal Kishore",
    "Sahaj Kulshrestha",
    "Parthivi Jain",
    "Varan Singh Rohila",
    "Achyut Sharma",
    "Priyanka Kumar",
    "Rishi Kumar"
)

$img = "../alumni/binod.jpg"
$url = "https://istenith.com/prody/"
$position = "Development Head, Byjus"

$startRow = 10
for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $startRow + $i
    $name = $names[$i]

    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 1).WrapText = $true

    $ws.Cells.Item($row, 2).Value = 2017
    $ws.Cells.Item($row, 3).Value = 2021

    $ws.Cells.Item($row, 4).Value = $img

    $ws.Cells.Item($row, 5).Value = $url
    $ws.Cells.Item($row, 6).Value = $url
    $ws.Cells.Item($row, 7).Value = $url

    $ws.Cells.Item($row, 8).Value = $position
}

# Hyperlinks were only added for the first four new rows (10-13) in the
# source edit.
for ($row = 10; $row -le 13; $row++) {
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $url, "", "", $url)
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $url, "", "", $url)
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 7), $url, "", "", $url)
}
